# Generate Report for Handback
# Updates the zh-cn and de-de localization-status sheets to reflect a
# handback: status flips to "Handed back: in sync with en-US", the
# "Latest Target File" / "Latest Handback File" columns (E/F) are filled
# in (with hyperlinks, mirroring columns A/C), and the handback datetime
# (column G) is stamped.

$wb = $excel.ActiveWorkbook

$sheetConfigs = @(
    @{
        SheetName   = "zh-cn"
        Row2Handoff = "487f546b-79dc-4d6e-9685-1fa141cba8f1.95e212ddbc3762be01df4d52572fa92e86bb7b16.zh-cn.xlf"
        Row3Handoff = "cd2cbb82-99d4-4be8-83d2-a81fa423dc03.96fb3de54468e89b7de6783c033bbe4bef4ff415.zh-cn.xlf"
        HandoffRepo = "OpenLocalizationTestOrg/olhandoff"
        HandoffRef  = "11139faa2f09ef73fd2b2210acac567b7edd127c"
        HandoffPath = "ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
        HandbackRepo = "OpenLocalizationTestOrg/olhandback"
        HandbackRef  = "11139faa2f09ef73fd2b2210acac567b7edd127c"
        HandbackPath = "ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht"
        HandbackDateTime = "2016-03-01 09:17:40"
    },
    @{
        SheetName   = "de-de"
        Row2Handoff = "487f546b-79dc-4d6e-9685-1fa141cba8f1.95e212ddbc3762be01df4d52572fa92e86bb7b16.de-de.xlf"
        Row3Handoff = "cd2cbb82-99d4-4be8-83d2-a81fa423dc03.96fb3de54468e89b7de6783c033bbe4bef4ff415.de-de.xlf"
        HandoffRepo = "OpenLocalizationTestOrg/olhandoff"
        HandoffRef  = "79eb258384feb2d48dc7d15ec8b8eb1697b25efc"
        HandoffPath = "ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"
        HandbackRepo = "OpenLocalizationTestOrg/olhandback"
        HandbackRef  = "79eb258384feb2d48dc7d15ec8b8eb1697b25efc"
        HandbackPath = "ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht"
        HandbackDateTime = "2016-03-01 09:17:59"
    }
)

$statusHandedBack = "Handed back: in sync with en-US"
$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/4112a2f41b163b04721de54f0fce106e44ee0c79/e2e"

foreach ($cfg in $sheetConfigs) {
    $ws = $wb.Worksheets.Item($cfg.SheetName)

    # Row 2 -> 487f546b-...md
    $mdName2 = $ws.Range("A2").Text
    $xlfName2 = $ws.Range("C2").Text

    $ws.Range("B2").Value = $statusHandedBack
    $ws.Range("G2").Value = $cfg.HandbackDateTime

    $ws.Range("E2").Value = $mdName2
    $ws.Range("F2").Value = $xlfName2
    $ws.Range("E2").Style = "HyperLink"
    $ws.Range("F2").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("E2"), "$mdBase/$mdName2", [Type]::Missing, [Type]::Missing, $mdName2) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/$($cfg.HandbackRepo)/blob/$($cfg.HandbackRef)/$($cfg.HandbackPath)/$xlfName2", [Type]::Missing, [Type]::Missing, $xlfName2) | Out-Null

    # Row 3 -> cd2cbb82-...md
    $mdName3 = $ws.Range("A3").Text
    $xlfName3 = $ws.Range("C3").Text

    $ws.Range("B3").Value = $statusHandedBack
    $ws.Range("G3").Value = $cfg.HandbackDateTime

    $ws.Range("E3").Value = $mdName3
    $ws.Range("F3").Value = $xlfName3
    $ws.Range("E3").Style = "HyperLink"
    $ws.Range("F3").Style = "HyperLink"

    $ws.Hyperlinks.Add($ws.Range("E3"), "$mdBase/$mdName3", [Type]::Missing, [Type]::Missing, $mdName3) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), "https://github.com/$($cfg.HandbackRepo)/blob/$($cfg.HandbackRef)/$($cfg.HandbackPath)/$xlfName3", [Type]::Missing, [Type]::Missing, $xlfName3) | Out-Null

    Write-Host "Updated sheet" $cfg.SheetName
}
